# Test Suite TS_9_1: split the first "Passed" test-result cell's run into
# "Passe" + "d" (two runs with identical formatting), matching the
# author's edit described in the commit message / diff.

$d = $word.ActiveDocument

# Locate the results table (Test Case ID / Incoming data / Expected Result / Test Result)
# and the first data row (TC-01), whose last cell holds the "Passed" verdict.
$table = $d.Tables.Item(2)
$cell = $table.Cell(2, 4)
$cellRange = $cell.Range

$start = $cellRange.Start
$target = $d.Range($start, $start + 5)   # covers "Passe" of "Passed"

if ($target.Text -ne "Passe") {
    throw "Unexpected cell contents, aborting: [" + $cellRange.Text + "]"
}

# Toggling a character property and reverting it forces the run to be
# split at this boundary while keeping formatting identical on both
# halves, producing two <w:r> runs: "Passe" and "d".
$target.Font.Bold = 1
$target.Font.Bold = 0

Write-Output ("Cell text after edit: [" + $cell.Range.Text + "]")
